# Reserva para multiples pasajeros
# Adds "part 2/3/4" columns (C, D, E) for the passenger first name, last
# name, and meal-option rows so the datapool sheet can drive multiple
# passengers in one booking row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: First name (Nombre) additional passengers ---------------------
$ws.Range("C2").Value = "Nombre p2"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("D2").Value = "Nombre p3"
$ws.Range("D2").Font.Name = "Ubuntu"
$ws.Range("D2").Font.Size = 13
$ws.Range("D2").Font.Bold = $false
$ws.Range("D2").Font.Underline = $false
$ws.Range("D2").Font.Color = 0
$ws.Range("D2").NumberFormat = "@"

$ws.Range("E2").Value = "Nombre p4"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- Row 3: Last name (Apellido) additional passengers ---------------------
$ws.Range("C3").Value = "Apellido p2"
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("D3").Value = "Apellido p3"
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("E3").Value = "Apellido p4"
$ws.Range("D2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- Row 4: Meal options for additional passengers --------------------------
$ws.Range("C4").Value = "Bland"
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("D4").Value = "Low Calorie"
$ws.Range("B4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("E4").Value = "Diabetic"
$ws.Range("B4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column width tweak (room for the 3 new narrower columns) --------------
# (Excel's ColumnWidth setter snaps to its internal pixel grid, so 10.8 is the
# closest character-width input that lands on the target ~11.6 stored width.)
$ws.Columns("C").ColumnWidth = 10.8

# --- View / selection -------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("E4").Select()

Write-Output "Reserva para multiples pasajeros: columns C-E added for rows 2-4"
